$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '38.582.73'
$ws.Range('E2').Value = '  -5.15%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.172.77'
$ws.Range('E3').Value = '  -8.46%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '291.42'
$ws.Range('E5').Value = '  -6.25%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '79.26'
$ws.Range('E6').Value = '  -9.24%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.498'
$ws.Range('E7').Value = '  -5.65%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -8.97%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0762'
$ws.Range('E10').Value = '  -9.05%  '
$ws.Range('B11').Value = 'OKB'
$ws.Range('C11').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '46.35'
$ws.Range('E11').Value = '  -12.04%  '
$ws.Range('B12').Value = 'Avalanche'
$ws.Range('C12').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '27.50'
$ws.Range('E12').Value = '  -10.50%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.106'
$ws.Range('E13').Value = '  -2.72%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.512.26'
$ws.Range('E14').Value = '  -8.41%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.02'
$ws.Range('E15').Value = '  -8.35%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '13.61'
$ws.Range('E16').Value = '  -9.39%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.194.21'
$ws.Range('E17').Value = '  -6.29%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.694'
$ws.Range('E18').Value = '  -9.00%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '38.490.85'
$ws.Range('E19').Value = '  -5.25%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0848'
$ws.Range('E20').Value = '  -6.95%  '
$ws.Range('E21').Value = '  -9.42%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '63.35'
$ws.Range('E22').Value = '  -8.02%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.59'
$ws.Range('E23').Value = '  -12.17%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '219.66'
$ws.Range('E24').Value = '  -5.68%  '
$ws.Range('E25').Value = '  -0.01%  '
$ws.Range('E26').Value = '  -11.71%  '
$ws.Range('E27').Value = '  -5.81%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '21.89'
$ws.Range('E28').Value = '  -8.19%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.16'
$ws.Range('E29').Value = '  -2.34%  '
$ws.Range('E30').Value = '  -6.82%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '146.60'
$ws.Range('E31').Value = '  -4.03%  '
$ws.Range('B32').Value = 'InjectiveProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '30.67'
$ws.Range('E32').Value = '  -9.33%  '
$ws.Range('B33').Value = 'FirstDigitalUSD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  -0.24%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.65'
$ws.Range('E34').Value = '  -11.63%  '
$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.30'
$ws.Range('E35').Value = '  -6.00%  '
$ws.Range('E36').Value = '  -7.93%  '
$ws.Range('E37').Value = '  -5.46%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0935'
$ws.Range('E38').Value = '  -6.19%  '
$ws.Range('E39').Value = '  -7.72%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.56'
$ws.Range('E40').Value = '  -9.07%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '13.91'
$ws.Range('E41').Value = '  -13.34%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.55'
$ws.Range('E42').Value = '  -8.10%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.867.95'
$ws.Range('E43').Value = '  -4.79%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.03'
$ws.Range('E44').Value = '  -14.57%  '
$ws.Range('E45').Value = '  -7.53%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '15.70'
$ws.Range('E46').Value = '  -11.31%  '
$ws.Range('E47').Value = '  -9.20%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.48'
$ws.Range('E48').Value = '  -9.23%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.389.36'
$ws.Range('E49').Value = '  -8.56%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '68.94'
$ws.Range('E50').Value = '  -4.98%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '84.62'
$ws.Range('E51').Value = '  -9.57%  '
